$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 247.5625
$ws.Range("I33").Value = 164.06667
$ws.Range("K33").Value = 164.06667
$ws.Range("M33").Value = 64.93333000000001

$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()

$ws.Range("H51").Value = 3500
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 3500
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 3500
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -4468

$ws.Range("H64").Value = 4633.3335
$ws.Range("I64").Value = 4150
$ws.Range("K64").Value = 4150
$ws.Range("M64").Value = -3902

$ws.Range("H67").Value = 4633.3335
$ws.Range("I67").Value = 4150
$ws.Range("K67").Value = 4150
$ws.Range("M67").Value = -3292

$ws.Range("H113").Value = 5772.579
$ws.Range("I113").Value = 4461.909
$ws.Range("K113").Value = 4461.909
$ws.Range("M113").Value = -1207.909

$ws.Range("H116").Value = 3236.182
$ws.Range("I116").Value = 2153.25
$ws.Range("J116").Value = 3855
$ws.Range("K116").Value = 2153.25
$ws.Range("L116").Value = 3855
$ws.Range("M116").Value = 1288.75
$ws.Range("N116").Value = -10739

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2373.25
$ws.Range("I32").Value = 2089
$ws.Range("K32").Value = 2089
$ws.Range("M32").Value = -1802

$ws.Range("H61").Value = 3440.125
$ws.Range("I61").Value = 2503
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 2503
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -2291
$ws.Range("N61").Value = -10424

$ws.Range("H74").Value = 2857.2222
$ws.Range("I74").Value = 2900.1428
$ws.Range("J74").Value = 2707
$ws.Range("K74").Value = 2900.1428
$ws.Range("L74").Value = 2707
$ws.Range("M74").Value = -2026.1428
$ws.Range("N74").Value = -4455

$ws.Range("H77").Value = 2857.2222
$ws.Range("I77").Value = 2900.1428
$ws.Range("J77").Value = 2707
$ws.Range("K77").Value = 14500.714
$ws.Range("L77").Value = 13535
$ws.Range("M77").Value = -10132.714
$ws.Range("N77").Value = -22271

$ws.Range("H101").Value = 21602
$ws.Range("J101").Value = 21602
$ws.Range("L101").Value = 21602
$ws.Range("N101").Value = -28092

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H122").Value = 3143.6667
$ws.Range("I122").Value = 3022.75
$ws.Range("J122").Value = 4111
$ws.Range("K122").Value = 9068.25
$ws.Range("L122").Value = 12333
$ws.Range("M122").Value = -6618.25
$ws.Range("N122").Value = -17233

$ws.Range("H136").Value = 3440.125
$ws.Range("I136").Value = 2503
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 7509
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -4959
$ws.Range("N136").Value = -35100

$ws.Range("H139").Value = 44999
$ws.Range("J139").Value = 44999
$ws.Range("L139").Value = 44999
$ws.Range("N139").Value = -55279

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 238.6
$ws.Range("I64").Value = 248.33333
$ws.Range("K64").Value = 248.33333
$ws.Range("M64").Value = -23.33332999999999

$ws.Range("H67").Value = 238.6
$ws.Range("I67").Value = 248.33333
$ws.Range("K67").Value = 248.33333
$ws.Range("M67").Value = 531.6666700000001

$ws.Range("H86").Value = 5286.5713
$ws.Range("J86").Value = 3250
$ws.Range("L86").Value = 3250
$ws.Range("N86").Value = -5496

$ws.Range("H89").Value = 5286.5713
$ws.Range("J89").Value = 3250
$ws.Range("L89").Value = 16250
$ws.Range("N89").Value = -27482

$ws.Range("H94").Value = 933.3333
$ws.Range("I94").Value = 933.3333
$ws.Range("K94").Value = 933.3333
$ws.Range("M94").Value = -482.3333

$ws.Range("H105").Value = 3817.1428
$ws.Range("I105").Value = 4070
$ws.Range("K105").Value = 4070
$ws.Range("M105").Value = -2323

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3002.5
$ws.Range("J62").Value = 2000
$ws.Range("L62").Value = 2000
$ws.Range("N62").Value = -3248

$ws.Range("H65").Value = 3002.5
$ws.Range("J65").Value = 2000
$ws.Range("L65").Value = 10000
$ws.Range("N65").Value = -16240

$ws.Range("H138").Value = 55000
$ws.Range("J138").Value = 55000
$ws.Range("L138").Value = 55000
$ws.Range("N138").Value = -65280

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1999
$ws.Range("I139").Value = 1998
$ws.Range("K139").Value = 5994
$ws.Range("M139").Value = -854

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 14081.363
$ws.Range("I46").Value = 5815.8335
$ws.Range("J46").Value = 24000
$ws.Range("K46").Value = 5815.8335
$ws.Range("L46").Value = 24000
$ws.Range("M46").Value = -5659.8335
$ws.Range("N46").Value = -24312

$ws.Range("H86").Value = 60000
$ws.Range("J86").Value = 60000
$ws.Range("L86").Value = 60000
$ws.Range("N86").Value = -62372

$ws.Range("H89").Value = 60000
$ws.Range("J89").Value = 60000
$ws.Range("L89").Value = 180000
$ws.Range("N89").Value = -191856

$ws.Range("H97").Value = 551.2222
$ws.Range("I97").Value = 561.25
$ws.Range("J97").Value = 471
$ws.Range("K97").Value = 561.25
$ws.Range("L97").Value = 471
$ws.Range("M97").Value = -65.25
$ws.Range("N97").Value = -1463

$ws.Range("H113").Value = 974.75
$ws.Range("I113").Value = 974.75
$ws.Range("K113").Value = 974.75
$ws.Range("M113").Value = 1195.25

$ws.Range("H122").Value = 1019.44446
$ws.Range("I122").Value = 1010.6667
$ws.Range("K122").Value = 3032.0001
$ws.Range("M122").Value = -582.0001000000002

$ws.Range("H126").Value = 10944.143
$ws.Range("I126").Value = 7527.5
$ws.Range("K126").Value = 22582.5
$ws.Range("M126").Value = -20112.5

$ws.Range("H132").Value = 1550
$ws.Range("I132").Value = 300
$ws.Range("J132").Value = 2800
$ws.Range("K132").Value = 900
$ws.Range("L132").Value = 8400
$ws.Range("M132").Value = 1630
$ws.Range("N132").Value = -13460

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4727.625
$ws.Range("I40").Value = 4376.1333
$ws.Range("K40").Value = 4376.1333
$ws.Range("M40").Value = -4240.1333

$ws.Range("H46").Value = 3539.4666
$ws.Range("J46").Value = 5000
$ws.Range("L46").Value = 5000
$ws.Range("N46").Value = -5376

$ws.Range("H61").Value = 2547.4
$ws.Range("I61").Value = 2547.4
$ws.Range("K61").Value = 2547.4
$ws.Range("M61").Value = -2345.4

$ws.Range("H100").Value = 2964.3333
$ws.Range("I100").Value = 2964.3333
$ws.Range("K100").Value = 2964.3333
$ws.Range("M100").Value = -2423.3333

$ws.Range("H113").Value = 2547.4
$ws.Range("I113").Value = 2547.4
$ws.Range("K113").Value = 2547.4
$ws.Range("M113").Value = -377.4000000000001

$ws.Range("H134").Value = 43429
$ws.Range("J134").Value = 43429
$ws.Range("L134").Value = 43429
$ws.Range("N134").Value = -53569

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 594.2
$ws.Range("I81").Value = 492.75
$ws.Range("J81").Value = 1000
$ws.Range("K81").Value = 985.5
$ws.Range("L81").Value = 2000
$ws.Range("M81").Value = 75.5
$ws.Range("N81").Value = -4122

$ws.Range("H84").Value = 594.2
$ws.Range("I84").Value = 492.75
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 4927.5
$ws.Range("L84").Value = 10000
$ws.Range("M84").Value = 376.5
$ws.Range("N84").Value = -20608

$ws.Range("H107").Value = 1859.04
$ws.Range("I107").Value = 1951.8235
$ws.Range("K107").Value = 5855.470499999999
$ws.Range("M107").Value = -3935.470499999999

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H132").Value = 2527.111
$ws.Range("I132").Value = 2374.3333
$ws.Range("K132").Value = 7122.999899999999
$ws.Range("M132").Value = -4592.999899999999
